# Remove unused fields and add distribution_name to entitlements and deliveries tables
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "model": the authorizations table field list.
# Remove the unused "description", "report_version" and "summary_form_id"
# fields, and add "distribution_name" after "distribution_id".
# ---------------------------------------------------------------------
$model = $wb.Worksheets.Item("model")

# Drop the trailing 4 rows (15-18); the remaining rows will be
# re-populated below with the final field list.
$model.Range("A15:B18").EntireRow.Delete()

# Row 3 used to be "name" (optional-field style) and must become
# "custom_delivery_form_id" using the same style as the other
# required fields (copy B2's format onto B3).
$model.Range("B2").Copy()
$model.Range("B3").PasteSpecial(-4122)

# Rows 10-12 move from the "required" style block into the "optional"
# style block (matching row 13's look).
$model.Range("B13").Copy()
$model.Range("B10:B12").PasteSpecial(-4122)

$model.Range("B3").Value = "custom_delivery_form_id"
$model.Range("B4").Value = "item_pack_id"
$model.Range("B5").Value = "item_pack_name"
$model.Range("B6").Value = "item_pack_description"
$model.Range("B7").Value = "item_pack_ranges"
$model.Range("B8").Value = "extra_field_entitlements"
$model.Range("B9").Value = "status"
$model.Range("B10").Value = "status_reason"
$model.Range("B11").Value = "date_created"
$model.Range("B12").Value = "for_member"
$model.Range("B13").Value = "distribution_id"
$model.Range("B14").Value = "distribution_name"

# ---------------------------------------------------------------------
# Sheet "settings"
# ---------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Range("B6").Value = "item_pack_name"

# ---------------------------------------------------------------------
# Selection / active-tab bookkeeping to match the edited file.
# "survey" ends up being the active/selected tab, so it must be the
# last sheet touched; "properties" keeps its original selection (E4)
# and is left alone.
# ---------------------------------------------------------------------
$model.Range("H8").Select()
$settings.Range("B7").Select()

$survey = $wb.Worksheets.Item("survey")
$survey.Range("F32").Select()
